$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.212.42"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.864.74"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7151"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07749"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3070"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08248"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "1.876.05"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7158"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.209"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "29.218.75"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.830"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007776"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("D20").Value = "2.123.56"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.943"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1587"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.899"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.493"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.303"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.341"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.077"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05186"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.913"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.54%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.175"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7273"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.689"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "1.165.02"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9026"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.082"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("D46").Value = "2.018.26"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5285"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.759"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.243"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.865"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9979"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.06%  "
